$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values to row 6 (X6, Y6)
$ws.Range("X6").Value = -0.35999999999999943
$ws.Range("Y6").Value = "Down"

# Add new row 7 of data
$ws.Range("A7").Value = 42648.885277777779
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("B7").Value = -4
$ws.Range("C7").Value = "Neutral"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = "Random"
$ws.Range("Q7").Value = 52.976913006825477
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = -0.0944
$ws.Range("S7").NumberFormat = "0.00%"
$ws.Range("T7").Value = -0.0257
$ws.Range("T7").NumberFormat = "0.00%"
$ws.Range("U7").Value = 6.62
$ws.Range("V7").Value = 1.88
$ws.Range("W7").Value = -2
